# Fix: capital costs page error when Sandia PV module model selected.
# - Mark "PV self-shading testing" (row 34) as Done.
# - Insert a new TODO row for re-arranging self-shading inputs in the UI.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Mark the "PV self-shading testing" task as Done.
$ws.Range("A34").Value() = "Done"

# Insert a new row above the old row 35 (shifts rows 35:61 down to 36:62,
# carrying formatting down from the row above and auto-adjusting the
# SUM(D17:D39) formula in H17 to SUM(D17:D40)).
$ws.Rows.Item(35).Insert()

# Populate the newly inserted row 35 with the new TODO item.
$ws.Range("A35").Value() = "Not done"
$ws.Range("B35").Value() = "Re-arrange self-shading inputs in UI with system design? Check inputs for usability in SDK"
$ws.Range("C35").Value() = "Janine"

# Update the view to match what the editor had selected/scrolled to.
$ws.Application.ActiveWindow.ScrollRow = 28
$ws.Range("B36").Select()
